$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H57").Value = 10733
$ws.Range("J57").Value = 11643.75
$ws.Range("L57").Value = 34931.25
$ws.Range("N57").Value = -35929.25

$ws.Range("H132").Value = 38285.56
$ws.Range("I132").Value = 55667.78
$ws.Range("J132").Value = 1940.909
$ws.Range("K132").Value = 167003.34
$ws.Range("L132").Value = 5822.727000000001
$ws.Range("M132").Value = -164473.34
$ws.Range("N132").Value = -10882.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3399.861
$ws.Range("I61").Value = 3528.3928
$ws.Range("J61").Value = 2950
$ws.Range("K61").Value = 3528.3928
$ws.Range("L61").Value = 2950
$ws.Range("M61").Value = -3316.3928
$ws.Range("N61").Value = -3374

$ws.Range("H96").Value = 9922
$ws.Range("J96").Value = 9922
$ws.Range("L96").Value = 9922
$ws.Range("N96").Value = -15414

$ws.Range("H132").Value = 1704.8909
$ws.Range("I132").Value = 1296.075
$ws.Range("J132").Value = 2795.0667
$ws.Range("K132").Value = 3888.225
$ws.Range("L132").Value = 8385.2001
$ws.Range("M132").Value = -1358.225
$ws.Range("N132").Value = -13445.2001

$ws.Range("H136").Value = 3399.861
$ws.Range("I136").Value = 3528.3928
$ws.Range("J136").Value = 2950
$ws.Range("K136").Value = 10585.1784
$ws.Range("L136").Value = 8850
$ws.Range("M136").Value = -8035.178400000001
$ws.Range("N136").Value = -13950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1245.6111
$ws.Range("I134").Value = 993.62964
$ws.Range("J134").Value = 2001.5555
$ws.Range("K134").Value = 2980.88892
$ws.Range("L134").Value = 6004.666499999999
$ws.Range("M134").Value = -445.8889199999999
$ws.Range("N134").Value = -11074.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1950.94
$ws.Range("I31").Value = 1277.92
$ws.Range("J31").Value = 2623.96
$ws.Range("K31").Value = 1277.92
$ws.Range("L31").Value = 2623.96
$ws.Range("M31").Value = -982.9200000000001
$ws.Range("N31").Value = -3213.96

$ws.Range("H34").Value = 1950.94
$ws.Range("I34").Value = 1277.92
$ws.Range("J34").Value = 2623.96
$ws.Range("K34").Value = 1277.92
$ws.Range("L34").Value = 2623.96
$ws.Range("M34").Value = -1075.92
$ws.Range("N34").Value = -3027.96

$ws.Range("H58").Value = 2414.3157
$ws.Range("I58").Value = 1847.1177
$ws.Range("J58").Value = 2873.476
$ws.Range("K58").Value = 1847.1177
$ws.Range("L58").Value = 2873.476
$ws.Range("M58").Value = -1644.1177
$ws.Range("N58").Value = -3279.476

$ws.Range("I86").Value = 35717876
$ws.Range("K86").Value = 35717876
$ws.Range("M86").Value = -35716753

$ws.Range("I89").Value = 35717876
$ws.Range("K89").Value = 178589380
$ws.Range("M89").Value = -178583764

$ws.Range("H99").Value = 7054.476
$ws.Range("I99").Value = 1988.8572
$ws.Range("J99").Value = 17185.715
$ws.Range("K99").Value = 1988.8572
$ws.Range("L99").Value = 17185.715
$ws.Range("M99").Value = -490.8571999999999
$ws.Range("N99").Value = -20181.715

$ws.Range("H126").Value = 7054.476
$ws.Range("I126").Value = 1988.8572
$ws.Range("J126").Value = 17185.715
$ws.Range("K126").Value = 5966.571599999999
$ws.Range("L126").Value = 51557.145
$ws.Range("M126").Value = -3496.571599999999
$ws.Range("N126").Value = -56497.145

$ws.Range("H132").Value = 835633.75
$ws.Range("I132").Value = 1390273.6
$ws.Range("K132").Value = 4170820.8
$ws.Range("M132").Value = -4168290.8

$ws.Range("H136").Value = 2414.3157
$ws.Range("I136").Value = 1847.1177
$ws.Range("J136").Value = 2873.476
$ws.Range("K136").Value = 5541.3531
$ws.Range("L136").Value = 8620.428
$ws.Range("M136").Value = -2991.3531
$ws.Range("N136").Value = -13720.428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2900
$ws.Range("J22").Value = 3250
$ws.Range("L22").Value = 9750
$ws.Range("N22").Value = -10088

$ws.Range("H27").Value = 2900
$ws.Range("J27").Value = 3250
$ws.Range("L27").Value = 9750
$ws.Range("N27").Value = -9954

$ws.Range("H109").Value = 5345
$ws.Range("I109").Value = 4633.3335
$ws.Range("J109").Value = 6199
$ws.Range("K109").Value = 13900.0005
$ws.Range("L109").Value = 18597
$ws.Range("M109").Value = -12860.0005
$ws.Range("N109").Value = -20677

$ws.Range("H113").Value = 674.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 674.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2023.5
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6363.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 12000
$ws.Range("J52").Value = 12000
$ws.Range("L52").Value = 12000
$ws.Range("N52").Value = -12518

$ws.Range("H54").Value = 13266.667
$ws.Range("J54").Value = 13266.667
$ws.Range("L54").Value = 13266.667
$ws.Range("N54").Value = -14046.667

$ws.Range("H70").Value = 5299.635
$ws.Range("I70").Value = 4776.35
$ws.Range("J70").Value = 5626.6875
$ws.Range("K70").Value = 4776.35
$ws.Range("L70").Value = 5626.6875
$ws.Range("M70").Value = -4506.35
$ws.Range("N70").Value = -6166.6875

$ws.Range("H73").Value = 5299.635
$ws.Range("I73").Value = 4776.35
$ws.Range("J73").Value = 5626.6875
$ws.Range("K73").Value = 4776.35
$ws.Range("L73").Value = 5626.6875
$ws.Range("M73").Value = -3840.35
$ws.Range("N73").Value = -7498.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58444.445
$ws.Range("I7").Value = 79523.08
$ws.Range("J7").Value = 3640
$ws.Range("K7").Value = 79523.08
$ws.Range("L7").Value = 3640
$ws.Range("M7").Value = -79411.08
$ws.Range("N7").Value = -3864

$ws.Range("H18").Value = 15000
$ws.Range("I18").Value = 15000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -14828
$ws.Range("N18").ClearContents()

$ws.Range("H22").Value = 663.4828
$ws.Range("I22").Value = 538.73334
$ws.Range("J22").Value = 797.1429000000001
$ws.Range("K22").Value = 538.73334
$ws.Range("L22").Value = 797.1429000000001
$ws.Range("M22").Value = -243.73334
$ws.Range("N22").Value = -1387.1429

$ws.Range("H27").Value = 663.4828
$ws.Range("I27").Value = 538.73334
$ws.Range("J27").Value = 797.1429000000001
$ws.Range("K27").Value = 538.73334
$ws.Range("L27").Value = 797.1429000000001
$ws.Range("M27").Value = -431.73334
$ws.Range("N27").Value = -1011.1429

$ws.Range("H45").Value = 21500
$ws.Range("I45").Value = 10000
$ws.Range("J45").Value = 33000
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 33000
$ws.Range("M45").Value = -9593
$ws.Range("N45").Value = -33814

$ws.Range("H46").Value = 1214.6154
$ws.Range("I46").Value = 1071.8182
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1071.8182
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -883.8181999999999
$ws.Range("N46").Value = -2376

$ws.Range("H55").Value = 195.35
$ws.Range("I55").Value = 188
$ws.Range("J55").Value = 204.33333
$ws.Range("K55").Value = 188
$ws.Range("L55").Value = 204.33333
$ws.Range("M55").Value = -15
$ws.Range("N55").Value = -550.3333299999999

$ws.Range("H68").Value = 1993
$ws.Range("J68").Value = 1993
$ws.Range("L68").Value = 1993
$ws.Range("N68").Value = -3491

$ws.Range("H71").Value = 1993
$ws.Range("J71").Value = 1993
$ws.Range("L71").Value = 9965
$ws.Range("N71").Value = -17453

$ws.Range("H122").Value = 27779620
$ws.Range("I122").Value = 55556560
$ws.Range("K122").Value = 166669680
$ws.Range("M122").Value = -166667230

$ws.Range("H126").Value = 58444.445
$ws.Range("I126").Value = 79523.08
$ws.Range("J126").Value = 3640
$ws.Range("K126").Value = 238569.24
$ws.Range("L126").Value = 10920
$ws.Range("M126").Value = -236099.24
$ws.Range("N126").Value = -15860

$ws.Range("H136").Value = 23813158
$ws.Range("I136").Value = 4036.9092
$ws.Range("J136").Value = 111113260
$ws.Range("K136").Value = 12110.7276
$ws.Range("L136").Value = 333339780
$ws.Range("M136").Value = -9560.7276
$ws.Range("N136").Value = -333344880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 600
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -826

$ws.Range("H46").Value = 38735.8
$ws.Range("J46").Value = 38735.8
$ws.Range("L46").Value = 38735.8
$ws.Range("N46").Value = -39197.8

$ws.Range("H62").Value = 2605.6155
$ws.Range("I62").Value = 2720
$ws.Range("J62").Value = 2534.125
$ws.Range("K62").Value = 2720
$ws.Range("L62").Value = 2534.125
$ws.Range("M62").Value = -2096
$ws.Range("N62").Value = -3782.125

$ws.Range("H65").Value = 2605.6155
$ws.Range("I65").Value = 2720
$ws.Range("J65").Value = 2534.125
$ws.Range("K65").Value = 13600
$ws.Range("L65").Value = 12670.625
$ws.Range("M65").Value = -10480
$ws.Range("N65").Value = -18910.625

$ws.Range("H122").Value = 73399.71000000001
$ws.Range("I122").Value = 112887.89
$ws.Range("J122").Value = 2321
$ws.Range("K122").Value = 338663.67
$ws.Range("L122").Value = 6963
$ws.Range("M122").Value = -336213.67
$ws.Range("N122").Value = -11863

$ws.Range("H126").Value = 67966.87
$ws.Range("I126").Value = 91863.45
$ws.Range("J126").Value = 2251.25
$ws.Range("K126").Value = 275590.35
$ws.Range("L126").Value = 6753.75
$ws.Range("M126").Value = -273120.35
$ws.Range("N126").Value = -11693.75

$ws.Range("H134").Value = 38735.8
$ws.Range("J134").Value = 38735.8
$ws.Range("L134").Value = 116207.4
$ws.Range("N134").Value = -121277.4

Write-Output "Applied all Carbuncle_Profits cell updates"